# Update cryptocurrency price (column D) and 1h volume change (column E) values
# to match the latest scrape, as captured by the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds prices formatted as plain text (e.g. "27.679.04", "0.3942").
# Several of the new values parse as plain numbers, so force those cells to
# the Text number format first -- otherwise Excel COM silently reinterprets
# the assigned string as a floating point number (losing trailing zeros /
# introducing binary rounding noise) instead of keeping it as text.
$textFormatRows = @(5, 8, 9, 11, 12, 14, 15, 16, 17, 18, 19, 20, 23, 24, 25, 27, 29, 30, 31, 32, 33, 34, 35, 36, 37, 38, 39, 40, 41, 42, 43, 44, 45, 46, 47, 49, 50, 51)
foreach ($r in $textFormatRows) {
    $ws.Range("D$r").NumberFormat = "@"
}

# --- Column D: Price ---
$ws.Range("D2").Value = "27.679.04"
$ws.Range("D3").Value = "1.874.06"
$ws.Range("D5").Value = "331.95"
$ws.Range("D8").Value = "0.3942"
$ws.Range("D9").Value = "47.94"
$ws.Range("D11").Value = "1.027"
$ws.Range("D12").Value = "21.99"
$ws.Range("D13").Value = "1.853.08"
$ws.Range("D14").Value = "5.962"
$ws.Range("D15").Value = "7.129"
$ws.Range("D16").Value = "1.006"
$ws.Range("D17").Value = "0.00001048"
$ws.Range("D18").Value = "87.07"
$ws.Range("D19").Value = "0.06664"
$ws.Range("D20").Value = "17.16"
$ws.Range("D22").Value = "27.691.92"
$ws.Range("D23").Value = "5.515"
$ws.Range("D24").Value = "10.99"
$ws.Range("D25").Value = "2.307"
$ws.Range("D26").Value = "2.090.05"
$ws.Range("D27").Value = "158.18"
$ws.Range("D29").Value = "2.101"
$ws.Range("D30").Value = "5.586"
$ws.Range("D31").Value = "122.17"
$ws.Range("D32").Value = "0.9746"
$ws.Range("D33").Value = "0.09544"
$ws.Range("D34").Value = "1.445"
$ws.Range("D35").Value = "3.592"
$ws.Range("D36").Value = "5.333"
$ws.Range("D37").Value = "0.06097"
$ws.Range("D38").Value = "0.02254"
$ws.Range("D39").Value = "1.227"
$ws.Range("D40").Value = "8.233"
$ws.Range("D41").Value = "0.6022"
$ws.Range("D42").Value = "0.1907"
$ws.Range("D43").Value = "10.23"
$ws.Range("D44").Value = "1.256"
$ws.Range("D45").Value = "0.5690"
$ws.Range("D46").Value = "12.24"
$ws.Range("D47").Value = "1.941"
$ws.Range("D49").Value = "115.45"
$ws.Range("D50").Value = "0.06885"
$ws.Range("D51").Value = "0.00000000303"

# --- Column E: Volume(1h) ---
$ws.Range("E2").Value = "  +1.11%  "
$ws.Range("E3").Value = "  +0.80%  "
$ws.Range("E4").Value = "  +0.28%  "
$ws.Range("E5").Value = "  +2.64%  "
$ws.Range("E6").Value = "  +0.22%  "
$ws.Range("E7").Value = "  +4.09%  "
$ws.Range("E8").Value = "  +2.04%  "
$ws.Range("E9").Value = "  -0.92%  "
$ws.Range("E10").Value = "  +1.46%  "
$ws.Range("E11").Value = "  +1.00%  "
$ws.Range("E12").Value = "  +2.88%  "
$ws.Range("E13").Value = "  -0.58%  "
$ws.Range("E14").Value = "  +0.80%  "
$ws.Range("E15").Value = "  +0.08%  "
$ws.Range("E16").Value = "  +0.36%  "
$ws.Range("E17").Value = "  +1.54%  "
$ws.Range("E18").Value = "  +1.45%  "
$ws.Range("E19").Value = "  +2.21%  "
$ws.Range("E20").Value = "  +0.56%  "
$ws.Range("E21").Value = "  +0.19%  "
$ws.Range("E22").Value = "  +1.18%  "
$ws.Range("E23").Value = "  -0.24%  "
$ws.Range("E24").Value = "  +1.35%  "
$ws.Range("E25").Value = "  +1.29%  "
$ws.Range("E26").Value = "  +0.30%  "
$ws.Range("E27").Value = "  +3.01%  "
$ws.Range("E28").Value = "  +2.34%  "
$ws.Range("E29").Value = "  +1.49%  "
$ws.Range("E30").Value = "  +2.66%  "
$ws.Range("E32").Value = "  +4.07%  "
$ws.Range("E33").Value = "  +2.59%  "
$ws.Range("E34").Value = "  -2.63%  "
$ws.Range("E35").Value = "  -0.07%  "
$ws.Range("E36").Value = "  +1.55%  "
$ws.Range("E37").Value = "  +1.73%  "
$ws.Range("E38").Value = "  +0.75%  "
$ws.Range("E39").Value = "  +0.46%  "
$ws.Range("E40").Value = "  -0.02%  "
$ws.Range("E41").Value = "  +1.97%  "
$ws.Range("E42").Value = "  +0.94%  "
$ws.Range("E43").Value = "  +1.23%  "
$ws.Range("E44").Value = "  -1.77%  "
$ws.Range("E45").Value = "  +1.47%  "
$ws.Range("E46").Value = "  +2.44%  "
$ws.Range("E47").Value = "  +0.92%  "
$ws.Range("E48").Value = "  +0.73%  "
$ws.Range("E49").Value = "  +6.67%  "
$ws.Range("E50").Value = "  +1.63%  "
$ws.Range("E51").Value = "  +13.11%  "

Write-Host "Updated cryptos list: 43 price cells and 49 volume cells."
